# Apply updated 2024 (column K) crime data for 2024-07-30 across all affected worksheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("K2").Value = 4581
$ws.Range("K3").Value = 4708
$ws.Range("J4").Value = 1830
$ws.Range("K4").Value = 959
$ws.Range("K5").Value = 340
$ws.Range("K6").Value = 5298
$ws.Range("J7").Value = 29301
$ws.Range("K7").Value = 15886

$ws = $wb.Worksheets.Item(2)
$ws.Range("K2").Value = 139
$ws.Range("K7").Value = 468
$ws.Range("K8").Value = 1069
$ws.Range("K10").Value = 86
$ws.Range("J11").Value = 540
$ws.Range("K11").Value = 313
$ws.Range("K12").Value = 27
$ws.Range("K13").Value = 17
$ws.Range("K15").Value = 162
$ws.Range("K18").Value = 108
$ws.Range("K19").Value = 477
$ws.Range("K20").Value = 368
$ws.Range("K23").Value = 164
$ws.Range("K25").Value = 78
$ws.Range("K29").Value = 842
$ws.Range("K31").Value = 179
$ws.Range("K33").Value = 667
$ws.Range("K36").Value = 201
$ws.Range("K42").Value = 588
$ws.Range("K51").Value = 199
$ws.Range("K52").Value = 419
$ws.Range("K54").Value = 298
$ws.Range("K55").Value = 183
$ws.Range("K56").Value = 18
$ws.Range("K63").Value = 48
$ws.Range("K65").Value = 365
$ws.Range("K67").Value = 612
$ws.Range("K72").Value = 73
$ws.Range("K77").Value = 115
$ws.Range("K79").Value = 394
$ws.Range("K83").Value = 339
$ws.Range("K85").Value = 722
$ws.Range("K89").Value = 227
$ws.Range("K90").Value = 143
$ws.Range("K91").Value = 171
$ws.Range("K92").Value = 60
$ws.Range("K93").Value = 60
$ws.Range("K94").Value = 198
$ws.Range("K95").Value = 278
$ws.Range("K97").Value = 129
$ws.Range("K99").Value = 267
$ws.Range("J101").Value = 29301
$ws.Range("K101").Value = 15886

$ws = $wb.Worksheets.Item(5)
$ws.Range("K2").Value = 164
$ws.Range("K3").Value = 148
$ws.Range("K7").Value = 468

$ws = $wb.Worksheets.Item(6)
$ws.Range("K2").Value = 103
$ws.Range("J4").Value = 32
$ws.Range("J7").Value = 540
$ws.Range("K7").Value = 313

$ws = $wb.Worksheets.Item(7)
$ws.Range("K2").Value = 58
$ws.Range("K3").Value = 72
$ws.Range("K7").Value = 227

$ws = $wb.Worksheets.Item(8)
$ws.Range("K2").Value = 252
$ws.Range("K3").Value = 242
$ws.Range("K4").Value = 41
$ws.Range("K6").Value = 165
$ws.Range("K7").Value = 722

$ws = $wb.Worksheets.Item(9)
$ws.Range("K2").Value = 112
$ws.Range("K7").Value = 419

$ws = $wb.Worksheets.Item(12)
$ws.Range("K2").Value = 296
$ws.Range("K3").Value = 322
$ws.Range("K6").Value = 363
$ws.Range("K7").Value = 1069

$ws = $wb.Worksheets.Item(13)
$ws.Range("K4").Value = 18
$ws.Range("K7").Value = 339

$ws = $wb.Worksheets.Item(14)
$ws.Range("K2").Value = 182
$ws.Range("K3").Value = 250
$ws.Range("K6").Value = 191
$ws.Range("K7").Value = 667

$ws = $wb.Worksheets.Item(15)
$ws.Range("K2").Value = 93
$ws.Range("K3").Value = 94
$ws.Range("K5").Value = 12
$ws.Range("K7").Value = 278

$ws = $wb.Worksheets.Item(17)
$ws.Range("K2").Value = 111
$ws.Range("K6").Value = 144
$ws.Range("K7").Value = 365

$ws = $wb.Worksheets.Item(18)
$ws.Range("K2").Value = 69
$ws.Range("K6").Value = 65
$ws.Range("K7").Value = 267

$ws = $wb.Worksheets.Item(20)
$ws.Range("K6").Value = 66
$ws.Range("K7").Value = 179

$ws = $wb.Worksheets.Item(21)
$ws.Range("K3").Value = 213
$ws.Range("K6").Value = 177
$ws.Range("K7").Value = 612

$ws = $wb.Worksheets.Item(24)
$ws.Range("K3").Value = 81
$ws.Range("K6").Value = 152
$ws.Range("K7").Value = 298

$ws = $wb.Worksheets.Item(25)
$ws.Range("K2").Value = 241
$ws.Range("K3").Value = 298
$ws.Range("K4").Value = 42
$ws.Range("K6").Value = 237
$ws.Range("K7").Value = 842

$ws = $wb.Worksheets.Item(27)
$ws.Range("K3").Value = 149
$ws.Range("K7").Value = 477

$ws = $wb.Worksheets.Item(32)
$ws.Range("K3").Value = 185
$ws.Range("K7").Value = 588

$ws = $wb.Worksheets.Item(33)
$ws.Range("K4").Value = 3
$ws.Range("K6").Value = 17

$ws = $wb.Worksheets.Item(34)
$ws.Range("K3").Value = 14
$ws.Range("K7").Value = 86

$ws = $wb.Worksheets.Item(36)
$ws.Range("K3").Value = 54
$ws.Range("K7").Value = 183

$ws = $wb.Worksheets.Item(39)
$ws.Range("K3").Value = 62
$ws.Range("K7").Value = 164

$ws = $wb.Worksheets.Item(40)
$ws.Range("K3").Value = 82
$ws.Range("K7").Value = 171

$ws = $wb.Worksheets.Item(42)
$ws.Range("K3").Value = 129
$ws.Range("K7").Value = 394

$ws = $wb.Worksheets.Item(44)
$ws.Range("K2").Value = 125
$ws.Range("K4").Value = 15
$ws.Range("K7").Value = 368

$ws = $wb.Worksheets.Item(45)
$ws.Range("K2").Value = 30
$ws.Range("K7").Value = 108

$ws = $wb.Worksheets.Item(47)
$ws.Range("K2").Value = 79
$ws.Range("K7").Value = 201

$ws = $wb.Worksheets.Item(48)
$ws.Range("K4").Value = 4
$ws.Range("K7").Value = 60

$ws = $wb.Worksheets.Item(51)
$ws.Range("K4").Value = 19
$ws.Range("K6").Value = 85
$ws.Range("K7").Value = 198

$ws = $wb.Worksheets.Item(52)
$ws.Range("K3").Value = 31
$ws.Range("K5").Value = 2
$ws.Range("K7").Value = 78

$ws = $wb.Worksheets.Item(54)
$ws.Range("K2").Value = 56
$ws.Range("K7").Value = 162

$ws = $wb.Worksheets.Item(64)
$ws.Range("K2").Value = 44
$ws.Range("K7").Value = 139

$ws = $wb.Worksheets.Item(65)
$ws.Range("K6").Value = 78
$ws.Range("K7").Value = 129

$ws = $wb.Worksheets.Item(66)
$ws.Range("K2").Value = 18
$ws.Range("K7").Value = 60

$ws = $wb.Worksheets.Item(74)
$ws.Range("K3").Value = 46
$ws.Range("K6").Value = 30
$ws.Range("K7").Value = 143

$ws = $wb.Worksheets.Item(75)
$ws.Range("K2").Value = 55
$ws.Range("K7").Value = 199

$ws = $wb.Worksheets.Item(82)
$ws.Range("K2").Value = 10
$ws.Range("K6").Value = 40
$ws.Range("K7").Value = 73

$ws = $wb.Worksheets.Item(84)
$ws.Range("K3").Value = 45
$ws.Range("K7").Value = 115

$ws = $wb.Worksheets.Item(86)
$ws.Range("K3").Value = 6
$ws.Range("K7").Value = 18

$ws = $wb.Worksheets.Item(91)
$ws.Range("K6").Value = 7
$ws.Range("K7").Value = 27
